# Auto commit at 2025-08-27  7:47:11.92
#
# Updates the monthly metric figures on the "Metrics" sheet (B2:B13), which
# ripple through the formulas on the "today" sheet (B11:B22, E11:E22,
# F11:F22 all reference Metrics!B2:B13). Also restores the selection/active
# sheet state captured at save time: "Metrics" has B2:B13 selected (but is
# no longer the active tab) while "today" becomes the active tab with D22
# selected.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 439509.50999999995
$metrics.Range("B3").Value = 377123.76
$metrics.Range("B4").Value = 139002.82000000004
$metrics.Range("B5").Value = 17262
$metrics.Range("B6").Value = 3835138.08
$metrics.Range("B7").Value = 3255838.4200000004
$metrics.Range("B8").Value = 1102645.3799999999
$metrics.Range("B9").Value = 147950
$metrics.Range("B10").Value = 32300461.879999995
$metrics.Range("B11").Value = 19285708.490000002
$metrics.Range("B12").Value = 11384354.270000001
$metrics.Range("B13").Value = 1245577

# Metrics keeps a selection of B2:B13 but relinquishes the active tab.
$metrics.Activate() | Out-Null
$metrics.Range("B2:B13").Select() | Out-Null

# "today" becomes the active sheet with D22 selected.
$today = $wb.Worksheets.Item("today")
$today.Activate() | Out-Null
$today.Range("D22").Select() | Out-Null
